# Update the "contacts" worksheet test data.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("contacts")

# Row 2
$ws.Range("A2").Value = "Mr."
$ws.Range("B2").Value = "Ram"
$ws.Range("C2").Value = "Selvan"
$ws.Range("D2").Value = "Salesforce"

# Row 3
$ws.Range("A3").Value = "Mr"
$ws.Range("B3").Value = "Chris"
$ws.Range("C3").Value = "David"
$ws.Range("D3").Value = "Walmart"

# Row 4
$ws.Range("A4").Value = "Mrs."
$ws.Range("B4").Value = "Anita"
$ws.Range("C4").Value = "Sharma"
$ws.Range("D4").Value = "Apple"

# Restore selection on the contacts sheet to A3, matching the saved state.
$ws.Activate()
$ws.Range("A3").Select()
